# chore: update Sheets via scheduled runner
# Refreshes cached market-price / profit figures (columns H-N) on a handful
# of leve rows across the per-job profit sheets. Values are static data
# (no formulas in this workbook), so cells are written/cleared directly.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 942
$ws.Range("I18").Value = 975
$ws.Range("J18").Value = 876
$ws.Range("K18").Value = 975
$ws.Range("L18").Value = 876
$ws.Range("M18").Value = -691
$ws.Range("N18").Value = -1444
$ws.Range("H19").Value = 364.8889
$ws.Range("I19").Value = 174
$ws.Range("J19").Value = 438.30768
$ws.Range("K19").Value = 174
$ws.Range("L19").Value = 438.30768
$ws.Range("M19").Value = 1
$ws.Range("N19").Value = -788.30768

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1233.3334
$ws.Range("I4").Value = 1100
$ws.Range("J4").Value = 1500
$ws.Range("K4").Value = 1100
$ws.Range("L4").Value = 1500
$ws.Range("M4").Value = -984
$ws.Range("N4").Value = -1732
$ws.Range("H5").Value = 581.8182
$ws.Range("I5").Value = 649.8
$ws.Range("J5").Value = 525.1667
$ws.Range("K5").Value = 649.8
$ws.Range("L5").Value = 525.1667
$ws.Range("M5").Value = -537.8
$ws.Range("N5").Value = -749.1667
$ws.Range("H37").Value = 29655
$ws.Range("J37").Value = 29655
$ws.Range("L37").Value = 29655
$ws.Range("N37").Value = -30201
$ws.Range("H44").Value = 36453.168
$ws.Range("J44").Value = 36453.168
$ws.Range("L44").Value = 36453.168
$ws.Range("N44").Value = -37429.168
$ws.Range("H55").Value = 49540
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 49540
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 49540
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -50170
$ws.Range("H63").Value = 2666.7
$ws.Range("I63").Value = 2382.2666
$ws.Range("J63").Value = 3520
$ws.Range("K63").Value = 2382.2666
$ws.Range("L63").Value = 3520
$ws.Range("M63").Value = -1696.2666
$ws.Range("N63").Value = -4892
$ws.Range("H66").Value = 2666.7
$ws.Range("I66").Value = 2382.2666
$ws.Range("J66").Value = 3520
$ws.Range("K66").Value = 11911.333
$ws.Range("L66").Value = 17600
$ws.Range("M66").Value = -8479.332999999999
$ws.Range("N66").Value = -24464
$ws.Range("H80").Value = 49996.668
$ws.Range("J80").Value = 49996.668
$ws.Range("L80").Value = 49996.668
$ws.Range("N80").Value = -51992.668
$ws.Range("H83").Value = 49996.668
$ws.Range("J83").Value = 49996.668
$ws.Range("L83").Value = 149990.004
$ws.Range("N83").Value = -159974.004
$ws.Range("H131").Value = 42804.668
$ws.Range("J131").Value = 42804.668
$ws.Range("L131").Value = 42804.668
$ws.Range("N131").Value = -52884.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 581.8182
$ws.Range("I4").Value = 649.8
$ws.Range("J4").Value = 525.1667
$ws.Range("K4").Value = 649.8
$ws.Range("L4").Value = 525.1667
$ws.Range("M4").Value = -534.8
$ws.Range("N4").Value = -755.1667
$ws.Range("H15").Value = 38666.668
$ws.Range("J15").Value = 38666.668
$ws.Range("L15").Value = 38666.668
$ws.Range("N15").Value = -39120.668
$ws.Range("H35").Value = 32270
$ws.Range("I35").Value = 30000
$ws.Range("J35").Value = 33026.668
$ws.Range("K35").Value = 30000
$ws.Range("L35").Value = 33026.668
$ws.Range("M35").Value = -29690
$ws.Range("N35").Value = -33646.668
$ws.Range("H130").Value = 50779.5
$ws.Range("J130").Value = 50779.5
$ws.Range("L130").Value = 50779.5
$ws.Range("N130").Value = -60819.5
$ws.Range("H132").Value = 50875
$ws.Range("J132").Value = 50875
$ws.Range("L132").Value = 50875
$ws.Range("N132").Value = -60995
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 1666.3334
$ws.Range("I35").Value = 1666.3334
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1666.3334
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -1372.3334
$ws.Range("N35").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 795.8333
$ws.Range("I68").Value = 800
$ws.Range("J68").Value = 792.8570999999999
$ws.Range("K68").Value = 2400
$ws.Range("L68").Value = 2378.5713
$ws.Range("M68").Value = -1589
$ws.Range("N68").Value = -4000.5713
$ws.Range("H71").Value = 795.8333
$ws.Range("I71").Value = 800
$ws.Range("J71").Value = 792.8570999999999
$ws.Range("K71").Value = 7200
$ws.Range("L71").Value = 7135.7139
$ws.Range("M71").Value = -3144
$ws.Range("N71").Value = -15247.7139
$ws.Range("H139").Value = 6407
$ws.Range("I139").Value = 6552.773
$ws.Range("J139").Value = 3200
$ws.Range("K139").Value = 19658.319
$ws.Range("L139").Value = 9600
$ws.Range("M139").Value = -14518.319
$ws.Range("N139").Value = -19880
$ws.Range("H141").Value = 34486416
$ws.Range("I141").Value = 45457816
$ws.Range("J141").Value = 4868.5713
$ws.Range("K141").Value = 136373448
$ws.Range("L141").Value = 14605.7139
$ws.Range("M141").Value = -136368268
$ws.Range("N141").Value = -24965.7139

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 4016117.2
$ws.Range("I2").Value = 4016117.2
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 4016117.2
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -4016004.2
$ws.Range("N2").ClearContents()
$ws.Range("H46").Value = 30198
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 30198
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 30198
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -30510
$ws.Range("H57").Value = 31277.143
$ws.Range("I57").Value = 18000
$ws.Range("J57").Value = 32298.46
$ws.Range("K57").Value = 18000
$ws.Range("L57").Value = 32298.46
$ws.Range("M57").Value = -17180
$ws.Range("N57").Value = -33938.46
$ws.Range("H80").Value = 389422.3
$ws.Range("I80").Value = 560443.3
$ws.Range("J80").Value = 4625
$ws.Range("K80").Value = 560443.3
$ws.Range("L80").Value = 4625
$ws.Range("M80").Value = -559445.3
$ws.Range("N80").Value = -6621
$ws.Range("H83").Value = 389422.3
$ws.Range("I83").Value = 560443.3
$ws.Range("J83").Value = 4625
$ws.Range("K83").Value = 2802216.5
$ws.Range("L83").Value = 23125
$ws.Range("M83").Value = -2797224.5
$ws.Range("N83").Value = -33109
$ws.Range("H126").Value = 7153.6
$ws.Range("I126").Value = 14740.25
$ws.Range("J126").Value = 2095.8333
$ws.Range("K126").Value = 44220.75
$ws.Range("L126").Value = 6287.499899999999
$ws.Range("M126").Value = -41750.75
$ws.Range("N126").Value = -11227.4999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 40000
$ws.Range("J20").Value = 40000
$ws.Range("L20").Value = 40000
$ws.Range("N20").Value = -40452
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H133").Value = 29940.445
$ws.Range("J133").Value = 29940.445
$ws.Range("L133").Value = 29940.445
$ws.Range("N133").Value = -35000.445
$ws.Range("H141").Value = 40000
$ws.Range("J141").Value = 40000
$ws.Range("L141").Value = 40000
$ws.Range("N141").Value = -50360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2451688.8
$ws.Range("I126").Value = 2674387.8
$ws.Range("K126").Value = 8023163.399999999
$ws.Range("M126").Value = -8020693.399999999
